$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.742.03'
$ws.Range("E2").Value = '  +7.60%  '
# Row 3
$ws.Range("D3").Value = '3.551.45'
$ws.Range("E3").Value = '  +10.13%  '
# Row 4
$ws.Range("E4").Value = '  -0.02%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '192.51'
$ws.Range("E5").Value = '  +10.36%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '560.86'
$ws.Range("E6").Value = '  +8.80%  '
# Row 7
$ws.Range("D7").Value = '3.546.74'
$ws.Range("E7").Value = '  +10.05%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.610'
$ws.Range("E8").Value = '  +3.22%  '
# Row 9
$ws.Range("E9").Value = '  -0.03%  '
# Row 10
$ws.Range("E10").Value = '  +7.33%  '
# Row 11
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.151'
$ws.Range("E11").Value = '  +16.40%  '
# Row 12
$ws.Range("B12").Value = 'Avalanche'
$ws.Range("C12").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '56.26'
$ws.Range("E12").Value = '  +6.92%  '
# Row 13
$ws.Range("E13").Value = '  +9.24%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.51'
$ws.Range("E14").Value = '  +7.25%  '
# Row 15
$ws.Range("D15").Value = '4.112.94'
$ws.Range("E15").Value = '  +10.01%  '
# Row 16
$ws.Range("D16").Value = '3.550.81'
$ws.Range("E16").Value = '  +10.25%  '
# Row 17
$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.122'
$ws.Range("E17").Value = '  +5.45%  '
# Row 18
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '67.767.10'
$ws.Range("E18").Value = '  +7.79%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '18.43'
$ws.Range("E19").Value = '  +7.58%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.95'
$ws.Range("E20").Value = '  +9.09%  '
# Row 21
$ws.Range("E21").Value = '  +4.81%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '407.28'
$ws.Range("E22").Value = '  +11.29%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.00'
$ws.Range("E23").Value = '  +7.93%  '
# Row 24
$ws.Range("E24").Value = '  +7.02%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.23'
$ws.Range("E25").Value = '  +8.70%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.41'
$ws.Range("E26").Value = '  +3.29%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.98'
$ws.Range("E27").Value = '  +14.55%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.15'
$ws.Range("E28").Value = '  +0.86%  '
# Row 29
$ws.Range("E29").Value = '  +7.04%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.85'
$ws.Range("E30").Value = '  +8.62%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '30.67'
$ws.Range("E31").Value = '  +8.94%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '692.18'
$ws.Range("E32").Value = '  +5.98%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.85'
$ws.Range("E33").Value = '  +8.91%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.83'
$ws.Range("E34").Value = '  +6.79%  '
# Row 35
$ws.Range("E35").Value = '  +8.36%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '60.87'
$ws.Range("E36").Value = '  +6.29%  '
# Row 37
$ws.Range("D37").Value = '0.0₃0835'
$ws.Range("E37").Value = '  +18.94%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '39.18'
$ws.Range("E38").Value = '  +7.31%  '
# Row 39
$ws.Range("E39").Value = '  -0.18%  '
# Row 40
$ws.Range("E40").Value = '  +6.95%  '
# Row 41
$ws.Range("E41").Value = '  +14.75%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.40'
$ws.Range("E42").Value = '  +18.90%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  +0.12%  '
# Row 44
$ws.Range("D44").Value = '3.064.22'
$ws.Range("E44").Value = '  +7.13%  '
# Row 45
$ws.Range("E45").Value = '  +15.88%  '
# Row 46
$ws.Range("E46").Value = '  +7.25%  '
# Row 47
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0422'
$ws.Range("E47").Value = '  +8.37%  '
# Row 48
$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.28'
$ws.Range("E48").Value = '  +12.00%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.13'
$ws.Range("E49").Value = '  +19.57%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.73'
$ws.Range("E50").Value = '  +2.39%  '
# Row 51
$ws.Range("E51").Value = '  +6.87%  '
